$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (row, date-serial, B, C, D)
$data = @(
    @(329, 44403, 3, 6, 33.37412392924686),
    @(330, 44404, 0, 6, 33.37412392924686),
    @(331, 44405, 2, 8, 44.49883190566248),
    @(332, 44406, 0, 8, 44.49883190566248),
    @(333, 44407, 2, 10, 55.6235398820781),
    @(334, 44408, 2, 11, 61.1858938702859),
    @(335, 44409, 1, 10, 55.6235398820781),
    @(336, 44410, 2, 9, 50.06118589387028),
    @(337, 44411, 1, 10, 55.6235398820781),
    @(338, 44412, 3, 11, 61.1858938702859),
    @(339, 44413, 0, 11, 61.1858938702859),
    @(340, 44414, 2, 11, 61.1858938702859),
    @(341, 44415, 2, 11, 61.1858938702859),
    @(342, 44416, 4, 14, 77.87295583490933),
    @(343, 44417, 6, 18, 100.1223717877406)
)

# Use row 328's column A as the style source (date column format)
$styleSource = $ws.Range("A328")

foreach ($row in $data) {
    $r = $row[0]
    $styleSource.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}

Write-Output "done"
